# Apply weekly price update edits for Fruta/Femacal de La Calera - Granada sheet.
# The underlying dataset rows (columns D, M-T) were re-shuffled/updated to a new
# weekly snapshot; apply the resulting target values to each affected cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "03/16/2021"
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 12000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 12000
$ws.Range("R2").Value = "Provincia del Elquí"
$ws.Range("S2").Value = 857
$ws.Range("D3").Value = "03/05/2021"
$ws.Range("M3").Value = 56
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 13000
$ws.Range("Q3").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R3").Value = "Provincia del Elquí"
$ws.Range("S3").Value = 929
$ws.Range("D4").Value = "02/11/2021"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "`$/caja 15 kilos granel"
$ws.Range("T4").Value = 15
$ws.Range("D6").Value = "04/26/2021"
$ws.Range("M6").Value = 68
$ws.Range("Q6").Value = "`$/caja 14 kilos granel"
$ws.Range("D7").Value = "05/07/2021"
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 14000
$ws.Range("Q7").Value = "`$/caja 14 kilos granel"
$ws.Range("R7").Value = "Provincia de Limarí"
$ws.Range("S7").Value = 1000
$ws.Range("D8").Value = "02/25/2021"
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 1000
$ws.Range("D10").Value = "02/15/2021"
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("Q10").Value = "`$/caja 15 kilos granel"
$ws.Range("S10").Value = 800
$ws.Range("T10").Value = 15
$ws.Range("D11").Value = "04/29/2021"
$ws.Range("M11").Value = 65
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 14000
$ws.Range("Q11").Value = "`$/caja 14 kilos granel"
$ws.Range("S11").Value = 1000
$ws.Range("T11").Value = 14
$ws.Range("D12").Value = "04/27/2021"
$ws.Range("M12").Value = 36
$ws.Range("N12").Value = 14000
$ws.Range("O12").Value = 14000
$ws.Range("P12").Value = 14000
$ws.Range("R12").Value = "Provincia de Limarí"
$ws.Range("S12").Value = 1000
$ws.Range("D13").Value = "05/06/2021"
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 14000
$ws.Range("O13").Value = 14000
$ws.Range("P13").Value = 14000
$ws.Range("Q13").Value = "`$/caja 14 kilos granel"
$ws.Range("S13").Value = 1000
$ws.Range("T13").Value = 14
$ws.Range("D14").Value = "03/04/2021"
$ws.Range("M14").Value = 80
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("Q14").Value = "`$/caja 15 kilos empedrada"
$ws.Range("S14").Value = 800
$ws.Range("D15").Value = "03/23/2021"
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = 13000
$ws.Range("O15").Value = 13000
$ws.Range("P15").Value = 13000
$ws.Range("Q15").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R15").Value = "Provincia del Elquí"
$ws.Range("S15").Value = 929
$ws.Range("D16").Value = "03/15/2021"
$ws.Range("M16").Value = 85
$ws.Range("N16").Value = 12000
$ws.Range("O16").Value = 12000
$ws.Range("P16").Value = 12000
$ws.Range("R16").Value = "Provincia del Elquí"
$ws.Range("S16").Value = 857
$ws.Range("D17").Value = "04/30/2021"
$ws.Range("M17").Value = 48
$ws.Range("D18").Value = "02/12/2021"
$ws.Range("M18").Value = 70
$ws.Range("D19").Value = "04/28/2021"
$ws.Range("M19").Value = 56
$ws.Range("N19").Value = 14000
$ws.Range("O19").Value = 14000
$ws.Range("P19").Value = 14000
$ws.Range("R19").Value = "Provincia de Limarí"
$ws.Range("S19").Value = 1000